$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for specific rows as per repull of data
$ws.Range("F2").Value = 12
$ws.Range("F4").Value = -6
$ws.Range("F5").Value = -7
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = -9
$ws.Range("F8").Value = 4
$ws.Range("F9").Value = 3
$ws.Range("F10").Value = 4
$ws.Range("F12").Value = -4
